# Case 3_26 (380 kV) simulation results update for pl_mw.xlsx
# Writes the refreshed power-flow results (columns B:D, F:G, J:M; rows 2-25)
# into the active worksheet, cell by cell, matching the new solved case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "B2"; Value = 1.687052136929452 },
    @{ Cell = "C2"; Value = 0.0142452460369995 },
    @{ Cell = "D2"; Value = 0.02906299646473087 },
    @{ Cell = "F2"; Value = 7.435570652922991 },
    @{ Cell = "G2"; Value = 0.002720002115830148 },
    @{ Cell = "J2"; Value = 0.2805919994846136 },
    @{ Cell = "K2"; Value = 1.139673038280108 },
    @{ Cell = "L2"; Value = 0.2463169711597715 },
    @{ Cell = "M2"; Value = 0.3657989241683524 },
    @{ Cell = "B3"; Value = 1.683526238653343 },
    @{ Cell = "C3"; Value = 0.01229971797557994 },
    @{ Cell = "D3"; Value = 0.02541701399940877 },
    @{ Cell = "F3"; Value = 7.240502031138448 },
    @{ Cell = "G3"; Value = 0.002725065492216183 },
    @{ Cell = "J3"; Value = 0.2774878871246287 },
    @{ Cell = "K3"; Value = 1.133179916801154 },
    @{ Cell = "L3"; Value = 0.2488286527273509 },
    @{ Cell = "M3"; Value = 0.3668667972428352 },
    @{ Cell = "B4"; Value = 1.682598546516772 },
    @{ Cell = "C4"; Value = 0.0111102631922293 },
    @{ Cell = "D4"; Value = 0.02317489508607196 },
    @{ Cell = "F4"; Value = 7.121119473286825 },
    @{ Cell = "G4"; Value = 0.002728336395837677 },
    @{ Cell = "J4"; Value = 0.2755796190430928 },
    @{ Cell = "K4"; Value = 1.130181604722068 },
    @{ Cell = "L4"; Value = 0.2505163392760039 },
    @{ Cell = "M4"; Value = 0.3677776925464578 },
    @{ Cell = "B5"; Value = 1.682531602769956 },
    @{ Cell = "C5"; Value = 0.01062672741637272 },
    @{ Cell = "D5"; Value = 0.02226020215401547 },
    @{ Cell = "F5"; Value = 7.072565207890818 },
    @{ Cell = "G5"; Value = 0.002729710183936132 },
    @{ Cell = "J5"; Value = 0.2748013339626496 },
    @{ Cell = "K5"; Value = 1.129208112596345 },
    @{ Cell = "L5"; Value = 0.2512407214084078 },
    @{ Cell = "M5"; Value = 0.3682130933944165 },
    @{ Cell = "B6"; Value = 1.682539276330942 },
    @{ Cell = "C6"; Value = 0.01054650460624629 },
    @{ Cell = "D6"; Value = 0.0221082529665253 },
    @{ Cell = "F6"; Value = 7.064508463789281 },
    @{ Cell = "G6"; Value = 0.002729940772997938 },
    @{ Cell = "J6"; Value = 0.274672059168239 },
    @{ Cell = "K6"; Value = 1.129061458183841 },
    @{ Cell = "L6"; Value = 0.2513632186697166 },
    @{ Cell = "M6"; Value = 0.3682892699626841 },
    @{ Cell = "B7"; Value = 1.682596384081791 },
    @{ Cell = "C7"; Value = 0.0111037374359384 },
    @{ Cell = "D7"; Value = 0.02316256349310208 },
    @{ Cell = "F7"; Value = 7.120464272515477 },
    @{ Cell = "G7"; Value = 0.002728354757312413 },
    @{ Cell = "J7"; Value = 0.2755691255340409 },
    @{ Cell = "K7"; Value = 1.13016747062089 },
    @{ Cell = "L7"; Value = 0.2505259601451968 },
    @{ Cell = "M7"; Value = 0.3677833045209589 },
    @{ Cell = "B8"; Value = 1.685579603012826 },
    @{ Cell = "C8"; Value = 0.01357330252067612 },
    @{ Cell = "D8"; Value = 0.02780649782023659 },
    @{ Cell = "F8"; Value = 7.368228062775358 },
    @{ Cell = "G8"; Value = 0.002721714435544965 },
    @{ Cell = "J8"; Value = 0.2795221456014261 },
    @{ Cell = "K8"; Value = 1.137228906247287 },
    @{ Cell = "L8"; Value = 0.2471528281345599 },
    @{ Cell = "M8"; Value = 0.3661141711567026 },
    @{ Cell = "B9"; Value = 1.701251215678781 },
    @{ Cell = "C9"; Value = 0.01846170618036069 },
    @{ Cell = "D9"; Value = 0.03689230374321539 },
    @{ Cell = "F9"; Value = 7.857352361602011 },
    @{ Cell = "G9"; Value = 0.002709971552365849 },
    @{ Cell = "J9"; Value = 0.2872587706116292 },
    @{ Cell = "K9"; Value = 1.158933218629784 },
    @{ Cell = "L9"; Value = 0.2416906065217717 },
    @{ Cell = "M9"; Value = 0.3648654548602259 },
    @{ Cell = "B10"; Value = 1.71876523322419 },
    @{ Cell = "C10"; Value = 0.02208834487259992 },
    @{ Cell = "D10"; Value = 0.04356454382467234 },
    @{ Cell = "F10"; Value = 8.218960621473627 },
    @{ Cell = "G10"; Value = 0.002702114660470775 },
    @{ Cell = "J10"; Value = 0.2929386021352869 },
    @{ Cell = "K10"; Value = 1.179693396760257 },
    @{ Cell = "L10"; Value = 0.2383774911905761 },
    @{ Cell = "M10"; Value = 0.3651823166050647 },
    @{ Cell = "B11"; Value = 1.72803903261925 },
    @{ Cell = "C11"; Value = 0.0237474354703977 },
    @{ Cell = "D11"; Value = 0.04660135896612871 },
    @{ Cell = "F11"; Value = 8.384012649488113 },
    @{ Cell = "G11"; Value = 0.002698705770776854 },
    @{ Cell = "J11"; Value = 0.2955226952257135 },
    @{ Cell = "K11"; Value = 1.190188849164912 },
    @{ Cell = "L11"; Value = 0.2370217400090482 },
    @{ Cell = "M11"; Value = 0.365594618688597 },
    @{ Cell = "B12"; Value = 1.73173886385888 },
    @{ Cell = "C12"; Value = 0.02437715499475246 },
    @{ Cell = "D12"; Value = 0.04775170870603063 },
    @{ Cell = "F12"; Value = 8.44659734941223 },
    @{ Cell = "G12"; Value = 0.00269743852979335 },
    @{ Cell = "J12"; Value = 0.2965013482389125 },
    @{ Cell = "K12"; Value = 1.194314816971655 },
    @{ Cell = "L12"; Value = 0.236530083130873 },
    @{ Cell = "M12"; Value = 0.3657893104535859 },
    @{ Cell = "B13"; Value = 1.7309336729611 },
    @{ Cell = "C13"; Value = 0.02424146688612439 },
    @{ Cell = "D13"; Value = 0.04750394134370595 },
    @{ Cell = "F13"; Value = 8.433114866134645 },
    @{ Cell = "G13"; Value = 0.002697710403881949 },
    @{ Cell = "J13"; Value = 0.2962905716769768 },
    @{ Cell = "K13"; Value = 1.193419469122915 },
    @{ Cell = "L13"; Value = 0.2366350040323937 },
    @{ Cell = "M13"; Value = 0.3657456650038391 },
    @{ Cell = "B14"; Value = 1.728339650299262 },
    @{ Cell = "C14"; Value = 0.02379921297236365 },
    @{ Cell = "D14"; Value = 0.04669599059906204 },
    @{ Cell = "F14"; Value = 8.38915984619689 },
    @{ Cell = "G14"; Value = 0.002698601041319259 },
    @{ Cell = "J14"; Value = 0.2956032069372299 },
    @{ Cell = "K14"; Value = 1.190525255279283 },
    @{ Cell = "L14"; Value = 0.2369808556612298 },
    @{ Cell = "M14"; Value = 0.3656098632445719 },
    @{ Cell = "B15"; Value = 1.726775230476761 },
    @{ Cell = "C15"; Value = 0.02352851308253889 },
    @{ Cell = "D15"; Value = 0.04620115047579532 },
    @{ Cell = "F15"; Value = 8.362247046043706 },
    @{ Cell = "G15"; Value = 0.002699149656138123 },
    @{ Cell = "J15"; Value = 0.2951821929544209 },
    @{ Cell = "K15"; Value = 1.188772215171213 },
    @{ Cell = "L15"; Value = 0.2371955295160149 },
    @{ Cell = "M15"; Value = 0.3655317027168508 },
    @{ Cell = "B16"; Value = 1.718185473724816 },
    @{ Cell = "C16"; Value = 0.02198011693548096 },
    @{ Cell = "D16"; Value = 0.04336612511698945 },
    @{ Cell = "F16"; Value = 8.208185526961756 },
    @{ Cell = "G16"; Value = 0.002702340753141414 },
    @{ Cell = "J16"; Value = 0.2927697376456138 },
    @{ Cell = "K16"; Value = 1.179028684202933 },
    @{ Cell = "L16"; Value = 0.2384691363140377 },
    @{ Cell = "M16"; Value = 0.3651607677675948 },
    @{ Cell = "B17"; Value = 1.71325069244449 },
    @{ Cell = "C17"; Value = 0.02103269566129029 },
    @{ Cell = "D17"; Value = 0.0416274375771053 },
    @{ Cell = "F17"; Value = 8.113818179394229 },
    @{ Cell = "G17"; Value = 0.00270434061867908 },
    @{ Cell = "J17"; Value = 0.2912899022629531 },
    @{ Cell = "K17"; Value = 1.173320918570283 },
    @{ Cell = "L17"; Value = 0.2392892053612385 },
    @{ Cell = "M17"; Value = 0.3650018931248056 },
    @{ Cell = "B18"; Value = 1.710535308288343 },
    @{ Cell = "C18"; Value = 0.02048863356736064 },
    @{ Cell = "D18"; Value = 0.04062753037274547 },
    @{ Cell = "F18"; Value = 8.059592578860588 },
    @{ Cell = "G18"; Value = 0.002705506450615674 },
    @{ Cell = "J18"; Value = 0.2904387650611255 },
    @{ Cell = "K18"; Value = 1.17013691676263 },
    @{ Cell = "L18"; Value = 0.2397751400289607 },
    @{ Cell = "M18"; Value = 0.3649357529147217 },
    @{ Cell = "B19"; Value = 1.7096370417282 },
    @{ Cell = "C19"; Value = 0.02030456924726565 },
    @{ Cell = "D19"; Value = 0.04028899854354506 },
    @{ Cell = "F19"; Value = 8.041241552928881 },
    @{ Cell = "G19"; Value = 0.00270590385791758 },
    @{ Cell = "J19"; Value = 0.2901505868392107 },
    @{ Cell = "K19"; Value = 1.169075853248131 },
    @{ Cell = "L19"; Value = 0.2399421181368879 },
    @{ Cell = "M19"; Value = 0.3649176940362473 },
    @{ Cell = "B20"; Value = 1.713763280175641 },
    @{ Cell = "C20"; Value = 0.0211334594210939 },
    @{ Cell = "D20"; Value = 0.04181250829401506 },
    @{ Cell = "F20"; Value = 8.123858338718321 },
    @{ Cell = "G20"; Value = 0.002704126120111836 },
    @{ Cell = "J20"; Value = 0.2914474301824299 },
    @{ Cell = "K20"; Value = 1.17391827608219 },
    @{ Cell = "L20"; Value = 0.2392004327753483 },
    @{ Cell = "M20"; Value = 0.3650161933375102 },
    @{ Cell = "B21"; Value = 1.729096472809942 },
    @{ Cell = "C21"; Value = 0.02392907315017112 },
    @{ Cell = "D21"; Value = 0.04693329400282664 },
    @{ Cell = "F21"; Value = 8.402068219588784 },
    @{ Cell = "G21"; Value = 0.002698338799149548 },
    @{ Cell = "J21"; Value = 0.295805099006067 },
    @{ Cell = "K21"; Value = 1.191371239619087 },
    @{ Cell = "L21"; Value = 0.2368786809976839 },
    @{ Cell = "M21"; Value = 0.3656487049334594 },
    @{ Cell = "B22"; Value = 1.74021375382938 },
    @{ Cell = "C22"; Value = 0.02576472206773417 },
    @{ Cell = "D22"; Value = 0.05028227052575573 },
    @{ Cell = "F22"; Value = 8.584380280139953 },
    @{ Cell = "G22"; Value = 0.00269469413137926 },
    @{ Cell = "J22"; Value = 0.2986537786448977 },
    @{ Cell = "K22"; Value = 1.203661389128285 },
    @{ Cell = "L22"; Value = 0.2354879646694386 },
    @{ Cell = "M22"; Value = 0.366286858540974 },
    @{ Cell = "B23"; Value = 1.734179896941754 },
    @{ Cell = "C23"; Value = 0.02478417955670409 },
    @{ Cell = "D23"; Value = 0.04849460450958532 },
    @{ Cell = "F23"; Value = 8.487031340203828 },
    @{ Cell = "G23"; Value = 0.002696626804491458 },
    @{ Cell = "J23"; Value = 0.2971332986538329 },
    @{ Cell = "K23"; Value = 1.197020931174848 },
    @{ Cell = "L23"; Value = 0.2362186358348382 },
    @{ Cell = "M23"; Value = 0.3659256964154523 },
    @{ Cell = "B24"; Value = 1.713531160231071 },
    @{ Cell = "C24"; Value = 0.02108790219450185 },
    @{ Cell = "D24"; Value = 0.04172883880947609 },
    @{ Cell = "F24"; Value = 8.1193190974343 },
    @{ Cell = "G24"; Value = 0.002704223044707911 },
    @{ Cell = "J24"; Value = 0.2913762129389994 },
    @{ Cell = "K24"; Value = 1.173647907173375 },
    @{ Cell = "L24"; Value = 0.2392405217985569 },
    @{ Cell = "M24"; Value = 0.3650096497192514 },
    @{ Cell = "B25"; Value = 1.695958916950417 },
    @{ Cell = "C25"; Value = 0.01713362825867648 },
    @{ Cell = "D25"; Value = 0.03443559714978051 },
    @{ Cell = "F25"; Value = 7.724655395515157 },
    @{ Cell = "G25"; Value = 0.002713012341840695 },
    @{ Cell = "J25"; Value = 0.285166922590065 },
    @{ Cell = "K25"; Value = 1.152218037724055 },
    @{ Cell = "L25"; Value = 0.243045174547099 },
    @{ Cell = "M25"; Value = 0.3649865303367683 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
